$d = $word.ActiveDocument

# --- Resize the header table's three columns ---
# 2970/5220/2880 (twips) -> 3780/3600/3690 (twips); COM Width is in points (twips/20)
$t = $d.Tables.Item(1)
$t.Columns.Item(1).Width = 189
$t.Columns.Item(2).Width = 180
$t.Columns.Item(3).Width = 184.5

# --- "... JUDICIAL CIRCUIT" -> "... JUDICIAL DISTRICT" ---
$d.Content.Find.Execute("CIRCUIT", $false, $false, $false, $false, $false, $true, 1, $false, "DISTRICT", 2) | Out-Null

# --- "ATTACHMENT TO ..." -> "ADDENDUM TO ..." ---
$d.Content.Find.Execute("ATTACHMENT", $false, $false, $false, $false, $false, $true, 1, $false, "ADDENDUM", 2) | Out-Null
